# Updated symbol list on Mon Jan  9 22:53:02 UTC 2023 with GitHub Actions
# Refresh Price / Volume(1h) figures in the crypto tracking sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "'0.79%"
$ws.Range("D3").Value = "'26.77"
$ws.Range("E3").Value = "'0.13%"
$ws.Range("D4").Value = "'4.897"
$ws.Range("E4").Value = "'3.61%"
$ws.Range("E5").Value = "'2.97%"
$ws.Range("D6").Value = "'6.909"
$ws.Range("E6").Value = "'2.47%"
$ws.Range("D7").Value = "'3.354"
$ws.Range("E7").Value = "'5.38%"
$ws.Range("D8").Value = "'1.341"
$ws.Range("E8").Value = "'48.95%"
$ws.Range("D9").Value = "'0.8836"
$ws.Range("E9").Value = "'3.30%"
$ws.Range("D10").Value = "'0.1463"
$ws.Range("E10").Value = "'2.47%"
$ws.Range("D11").Value = "'0.05102"
$ws.Range("E11").Value = "'0.56%"
$ws.Range("D12").Value = "'0.07397"
$ws.Range("E12").Value = "'3.43%"
$ws.Range("D13").Value = "'0.03137"
$ws.Range("E13").Value = "'-0.69%"
$ws.Range("D14").Value = "'0.09041"
$ws.Range("E14").Value = "'-0.13%"
$ws.Range("D15").Value = "'0.001557"
$ws.Range("E15").Value = "'1.24%"
$ws.Range("D16").Value = "'0.0006324"
$ws.Range("E16").Value = "'3.98%"
$ws.Range("D17").Value = "'0.006031"
$ws.Range("E17").Value = "'-0.19%"
$ws.Range("D18").Value = "'3.465"
$ws.Range("E18").Value = "'-0.01%"
$ws.Range("D21").Value = "'0.1333"
$ws.Range("E21").Value = "'4.02%"
$ws.Range("D22").Value = "'3.904"
$ws.Range("E22").Value = "'1.35%"
$ws.Range("D23").Value = "'0.04344"
$ws.Range("E23").Value = "'2.24%"
$ws.Range("D24").Value = "'0.001179"
$ws.Range("E24").Value = "'0.18%"
$ws.Range("D25").Value = "'0.003653"
$ws.Range("E25").Value = "'-11.94%"
$ws.Range("E27").Value = "'1.07%"
$ws.Range("D40").Value = "'0.04043"
$ws.Range("E40").Value = "'1.78%"
$ws.Range("D41").Value = "'0.006618"
$ws.Range("E41").Value = "'57.75%"
$ws.Range("D42").Value = "'0.1163"
$ws.Range("E42").Value = "'4.05%"
$ws.Range("E43").Value = "'4.43%"
$ws.Range("E44").Value = "'7.33%"
$ws.Range("D45").Value = "'0.00005335"
$ws.Range("E45").Value = "'3.55%"
$ws.Range("E46").Value = "'161.35%"
$ws.Range("D47").Value = "'0.02119"
$ws.Range("E47").Value = "'-29.24%"
